$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends a new pair of price rows (Primera/Segunda) for
# the latest sampling date, pushing all the existing history down by two
# rows. Insert two blank rows at 58:59 to make room, then populate them.
$ws.Rows("58:59").Insert()

# New "Primera" observation
$ws.Range("A58").Value = 11
$ws.Range("B58").Value = 'Vega Monumental Concepción'
$ws.Range("C58").Value = 'Bíobío'
$ws.Range("D58").Value = 44952
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 100112044
$ws.Range("G58").Value = 'Perejil'
$ws.Range("H58").Value = 'Sin especificar'
$ws.Range("I58").Value = 'Primera'
$ws.Range("J58").Value = 200
$ws.Range("K58").Value = 700
$ws.Range("L58").Value = 800
$ws.Range("M58").Value = 750
$ws.Range("N58").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O58").Value = 'Región de Ñuble'
$ws.Range("P58").Value = 750
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = 'Hortaliza'

# New "Segunda" observation
$ws.Range("A59").Value = 11
$ws.Range("B59").Value = 'Vega Monumental Concepción'
$ws.Range("C59").Value = 'Bíobío'
$ws.Range("D59").Value = 44952
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = 100112044
$ws.Range("G59").Value = 'Perejil'
$ws.Range("H59").Value = 'Sin especificar'
$ws.Range("I59").Value = 'Segunda'
$ws.Range("J59").Value = 100
$ws.Range("K59").Value = 600
$ws.Range("L59").Value = 600
$ws.Range("M59").Value = 600
$ws.Range("N59").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O59").Value = 'Región de Ñuble'
$ws.Range("P59").Value = 600
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = 'Hortaliza'
